# Update 'want to go' counts (column F) across sheets to match the newly generated site data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1858
$ws.Cells.Item(8, 6).Value = 3672
$ws.Cells.Item(9, 6).Value = 160
$ws.Cells.Item(11, 6).Value = 83
$ws.Cells.Item(12, 6).Value = 65
$ws.Cells.Item(13, 6).Value = 71
$ws.Cells.Item(14, 6).Value = 640
$ws.Cells.Item(15, 6).Value = 139
$ws.Cells.Item(16, 6).Value = 831
$ws.Cells.Item(17, 6).Value = 39
$ws.Cells.Item(19, 6).Value = 142
$ws.Cells.Item(20, 6).Value = 63
$ws.Cells.Item(22, 6).Value = 77
$ws.Cells.Item(23, 6).Value = 3022
$ws.Cells.Item(24, 6).Value = 5418
$ws.Cells.Item(27, 6).Value = 499
$ws.Cells.Item(28, 6).Value = 29
$ws.Cells.Item(29, 6).Value = 3148
$ws.Cells.Item(30, 6).Value = 325
$ws.Cells.Item(31, 6).Value = 2338
$ws.Cells.Item(33, 6).Value = 503
$ws.Cells.Item(34, 6).Value = 101
$ws.Cells.Item(35, 6).Value = 158
$ws.Cells.Item(36, 6).Value = 215
$ws.Cells.Item(37, 6).Value = 326
$ws.Cells.Item(38, 6).Value = 74
$ws.Cells.Item(39, 6).Value = 484
$ws.Cells.Item(40, 6).Value = 839
$ws.Cells.Item(45, 6).Value = 516

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 82

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1858
$ws.Cells.Item(8, 6).Value = 3672
$ws.Cells.Item(9, 6).Value = 160
$ws.Cells.Item(11, 6).Value = 83
$ws.Cells.Item(12, 6).Value = 82
$ws.Cells.Item(13, 6).Value = 65
$ws.Cells.Item(14, 6).Value = 71
$ws.Cells.Item(15, 6).Value = 640
$ws.Cells.Item(16, 6).Value = 139
$ws.Cells.Item(17, 6).Value = 831
$ws.Cells.Item(18, 6).Value = 39
$ws.Cells.Item(20, 6).Value = 142
$ws.Cells.Item(21, 6).Value = 63
$ws.Cells.Item(23, 6).Value = 77
$ws.Cells.Item(24, 6).Value = 3022
$ws.Cells.Item(25, 6).Value = 5418
$ws.Cells.Item(28, 6).Value = 499
$ws.Cells.Item(29, 6).Value = 29
$ws.Cells.Item(30, 6).Value = 3148
$ws.Cells.Item(31, 6).Value = 325
$ws.Cells.Item(32, 6).Value = 2338
$ws.Cells.Item(34, 6).Value = 503
$ws.Cells.Item(35, 6).Value = 101
$ws.Cells.Item(36, 6).Value = 158
$ws.Cells.Item(37, 6).Value = 215
$ws.Cells.Item(38, 6).Value = 326
$ws.Cells.Item(39, 6).Value = 74
$ws.Cells.Item(40, 6).Value = 484
$ws.Cells.Item(41, 6).Value = 839
$ws.Cells.Item(46, 6).Value = 516
